# Generate Report for Handoff
# Regenerate the localization-status report for the new source file GUID
# (198a50dc-ec52-4f16-8978-c18da6510493 -> 3d71b4c5-947a-429a-9c10-4e23a09dd6d6)
# and refreshed handoff/handback timestamps + xliff content hashes.

$wb = $excel.ActiveWorkbook

$oldGuid = "198a50dc-ec52-4f16-8978-c18da6510493"
$newGuid = "3d71b4c5-947a-429a-9c10-4e23a09dd6d6"

$oldZhXlf = "$oldGuid.1bc70fa953ac64a7f548d201d7dd17e97c54eeca.zh-cn.xlf"
$newZhXlf = "$newGuid.dc1311b846f9dd62cbf907a065a0b0c12964926f.zh-cn.xlf"

$oldDeXlf = "$oldGuid.1bc70fa953ac64a7f548d201d7dd17e97c54eeca.de-de.xlf"
$newDeXlf = "$newGuid.dc1311b846f9dd62cbf907a065a0b0c12964926f.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-21 07:03:36"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = "2016-08-21 07:03:32"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = "2016-08-21 07:03:36"

# ---------------------------------------------------------------------------
# Refresh the hyperlink display text on each sheet (keeps the same target
# address / relationship, only the visible text changes) by re-creating the
# hyperlink in place.
# ---------------------------------------------------------------------------
function Update-HyperlinkDisplay($ws, $cellAddr, $newDisplay) {
    $targetAddr = $null
    foreach ($h in $ws.Hyperlinks) {
        $targetAddr = $h.Address
    }
    foreach ($h in $ws.Hyperlinks) {
        $h.Delete()
    }
    if ($targetAddr) {
        $ws.Hyperlinks.Add($ws.Range($cellAddr), $targetAddr, "", "", $newDisplay)
    }
}

Update-HyperlinkDisplay $wsOverview "B2" "e2e\$newGuid.md"
Update-HyperlinkDisplay $wsZh "A2" "$newGuid.md"
Update-HyperlinkDisplay $wsDe "A2" "$newGuid.md"
